$wb = $excel.ActiveWorkbook

# Sheet "owhm_wb_dict" (2nd sheet) - rename the "Stream Leakage" name to "Stream Recharge"
# and add a new row for GHB_NET / black / Net Groundwater Flow.
$owhm = $wb.Worksheets.Item("owhm_wb_dict")

$owhm.Range("A15").Value = "GHB_NET"
$owhm.Range("B15").Value = "black"
$owhm.Range("C15").Value = "Net Groundwater Flow"

$owhm.Range("C4").Value = "Stream Recharge"

$owhm.Range("C5").Select()
$owhm.Activate()

$wb.Save()
